$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the crypto price refresh.
# Columns B/C hold plain text (coin name / link); no special handling needed.
# Columns D/E hold numeric-looking text (price / % change) that must stay as
# literal text (matching the original inlineStr values), so we force the
# cell to Text format before assigning the value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.31%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.45%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.018"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.50%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07843"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.12%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.165"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.79%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.037"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.01%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.054"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.17%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9231"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.40%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09930"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.29%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1871"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.84%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08688"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.13%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03583"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.88%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09935"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.04%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001490"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.31%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005625"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.49%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.50%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.14%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.87%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.913"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "7.99%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2201"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.65%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04601"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.75%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.005184"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "14.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001232"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.03%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.82%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002718"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.72%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01815"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.23%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04739"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.51%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007909"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.63%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1406"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.16%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007593"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.38%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002241"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.50%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01043"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "14.23%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006338"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.14%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005802"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.03%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.10"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "483.75%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.13%"
